$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.473752
$ws.Range("H2").Value = 10.421256
$ws.Range("I2").Value = 0.2912886159317403
$ws.Range("J2").Value = 0.2912886159317402
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.378907333333333
$ws.Range("N2").Value = 7.136722
$ws.Range("O2").Value = 0.3199597759432165
$ws.Range("P2").Value = 0.3199597759432165
$ws.Range("Q2").Value = 8.263734106981332
$ws.Range("R2").Value = 74.373606962832
$ws.Range("S2").Value = 0.09320064028832925
$ws.Range("T2").Value = 0.09320064028832922

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.473752
$ws.Range("H3").Value = 10.421256
$ws.Range("I3").Value = 0.2912886159317403
$ws.Range("J3").Value = 0.2912886159317402
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.056112666666666
$ws.Range("N3").Value = 15.168338
$ws.Range("O3").Value = 0.6800402240567835
$ws.Range("P3").Value = 0.6800402240567835
$ws.Range("Q3").Value = 17.56368148805867
$ws.Range("R3").Value = 158.073133392528
$ws.Range("S3").Value = 0.198087975643411
$ws.Range("T3").Value = 0.198087975643411

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.016808666666666
$ws.Range("H4").Value = 9.050426
$ws.Range("I4").Value = 0.2529720086650434
$ws.Range("J4").Value = 0.2529720086650434
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.378907333333333
$ws.Range("N4").Value = 7.136722
$ws.Range("O4").Value = 0.3199597759432165
$ws.Range("P4").Value = 0.3199597759432165
$ws.Range("Q4").Value = 7.176708260396889
$ws.Range("R4").Value = 64.590374343572
$ws.Range("S4").Value = 0.0809408672123727
$ws.Range("T4").Value = 0.0809408672123727

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.016808666666666
$ws.Range("H5").Value = 9.050426
$ws.Range("I5").Value = 0.2529720086650434
$ws.Range("J5").Value = 0.2529720086650434
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.056112666666666
$ws.Range("N5").Value = 15.168338
$ws.Range("O5").Value = 0.6800402240567835
$ws.Range("P5").Value = 0.6800402240567835
$ws.Range("Q5").Value = 15.25332451244311
$ws.Range("R5").Value = 137.279920611988
$ws.Range("S5").Value = 0.1720311414526707
$ws.Range("T5").Value = 0.1720311414526707

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.631845666666667
$ws.Range("H6").Value = 4.895537
$ws.Range("I6").Value = 0.136837075777874
$ws.Range("J6").Value = 0.136837075777874
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.378907333333333
$ws.Range("N6").Value = 7.136722
$ws.Range("O6").Value = 0.3199597759432165
$ws.Range("P6").Value = 0.3199597759432165
$ws.Range("Q6").Value = 3.882009623301556
$ws.Range("R6").Value = 34.938086609714
$ws.Range("S6").Value = 0.0437823601066135
$ws.Range("T6").Value = 0.0437823601066135

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.631845666666667
$ws.Range("H7").Value = 4.895537
$ws.Range("I7").Value = 0.136837075777874
$ws.Range("J7").Value = 0.136837075777874
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.056112666666666
$ws.Range("N7").Value = 15.168338
$ws.Range("O7").Value = 0.6800402240567835
$ws.Range("P7").Value = 0.6800402240567835
$ws.Range("Q7").Value = 8.250795545278445
$ws.Range("R7").Value = 74.25715990750601
$ws.Range("S7").Value = 0.0930547156712605
$ws.Range("T7").Value = 0.0930547156712605

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.673383
$ws.Range("H8").Value = 11.020149
$ws.Range("I8").Value = 0.3080285091904039
$ws.Range("J8").Value = 0.3080285091904039
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.378907333333333
$ws.Range("N8").Value = 7.136722
$ws.Range("O8").Value = 0.3199597759432165
$ws.Range("P8").Value = 0.3199597759432165
$ws.Range("Q8").Value = 8.738637756842
$ws.Range("R8").Value = 78.64773981157799
$ws.Range("S8").Value = 0.09855673278468462
$ws.Range("T8").Value = 0.09855673278468462

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.673383
$ws.Range("H9").Value = 11.020149
$ws.Range("I9").Value = 0.3080285091904039
$ws.Range("J9").Value = 0.3080285091904039
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.056112666666666
$ws.Range("N9").Value = 15.168338
$ws.Range("O9").Value = 0.6800402240567835
$ws.Range("P9").Value = 0.6800402240567835
$ws.Range("Q9").Value = 18.573038315818
$ws.Range("R9").Value = 167.157344842362
$ws.Range("S9").Value = 0.2094717764057193
$ws.Range("T9").Value = 0.2094717764057193

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.129675
$ws.Range("H10").Value = 0.389025
$ws.Range("I10").Value = 0.01087379043493848
$ws.Range("J10").Value = 0.01087379043493848
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.378907333333333
$ws.Range("N10").Value = 7.136722
$ws.Range("O10").Value = 0.3199597759432165
$ws.Range("P10").Value = 0.3199597759432165
$ws.Range("Q10").Value = 0.30848480845
$ws.Range("R10").Value = 2.77636327605
$ws.Range("S10").Value = 0.003479175551216407
$ws.Range("T10").Value = 0.003479175551216406

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.129675
$ws.Range("H11").Value = 0.389025
$ws.Range("I11").Value = 0.01087379043493848
$ws.Range("J11").Value = 0.01087379043493848
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.056112666666666
$ws.Range("N11").Value = 15.168338
$ws.Range("O11").Value = 0.6800402240567835
$ws.Range("P11").Value = 0.6800402240567835
$ws.Range("Q11").Value = 0.65565141005
$ws.Range("R11").Value = 5.90086269045
$ws.Range("S11").Value = 0.007394614883722075
$ws.Range("T11").Value = 0.007394614883722075

